$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text so numeric-looking price strings are not
# auto-converted to floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '26.361.73'
$ws.Range('E2').Value = '  -0.53%  '
$ws.Range('D3').Value = '1.829.94'
$ws.Range('E3').Value = '  -0.51%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = '251.24'
$ws.Range('E5').Value = '  -3.51%  '
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('D7').Value = '0.5243'
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').Value = '0.2761'
$ws.Range('E8').Value = '  -13.98%  '
$ws.Range('D9').Value = '0.06808'
$ws.Range('E9').Value = '  +0.30%  '
$ws.Range('D10').Value = '1.848.00'
$ws.Range('E10').Value = '  +0.45%  '
$ws.Range('D11').Value = '16.39'
$ws.Range('E11').Value = '  -12.65%  '
$ws.Range('D12').Value = '0.07078'
$ws.Range('E12').Value = '  -8.42%  '
$ws.Range('D13').Value = '0.6846'
$ws.Range('E13').Value = '  -12.78%  '
$ws.Range('D14').Value = '85.83'
$ws.Range('E14').Value = '  -2.12%  '
$ws.Range('D15').Value = '4.835'
$ws.Range('E15').Value = '  -3.49%  '
$ws.Range('D16').Value = '1.003'
$ws.Range('E16').Value = '  +0.26%  '
$ws.Range('D17').Value = '1.001'
$ws.Range('D18').Value = '13.13'
$ws.Range('E18').Value = '  -5.12%  '
$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = '26.397.19'
$ws.Range('E19').Value = '  -0.50%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.000007287'
$ws.Range('E20').Value = '  -8.08%  '
$ws.Range('D21').Value = '2.084.85'
$ws.Range('E21').Value = '  +0.41%  '
$ws.Range('D22').Value = '4.488'
$ws.Range('E22').Value = '  -2.87%  '
$ws.Range('D23').Value = '5.789'
$ws.Range('E23').Value = '  -3.15%  '
$ws.Range('D24').Value = '8.937'
$ws.Range('E24').Value = '  -4.67%  '
$ws.Range('D25').Value = '142.44'
$ws.Range('E25').Value = '  +0.85%  '
$ws.Range('D26').Value = '1.671'
$ws.Range('E26').Value = '  -0.86%  '
$ws.Range('D27').Value = '2.011'
$ws.Range('E27').Value = '  -6.81%  '
$ws.Range('D28').Value = '16.49'
$ws.Range('E28').Value = '  -2.54%  '
$ws.Range('D29').Value = '108.82'
$ws.Range('E29').Value = '  -2.61%  '
$ws.Range('D30').Value = '4.043'
$ws.Range('E30').Value = '  -2.48%  '
$ws.Range('D31').Value = '0.08709'
$ws.Range('E31').Value = '  +0.23%  '
$ws.Range('D32').Value = '3.835'
$ws.Range('E32').Value = '  -5.79%  '
$ws.Range('D33').Value = '0.04670'
$ws.Range('E33').Value = '  -3.89%  '
$ws.Range('D34').Value = '2.878'
$ws.Range('E34').Value = '  +0.84%  '
$ws.Range('D35').Value = '1.101'
$ws.Range('E35').Value = '  -3.02%  '
$ws.Range('D36').Value = '0.6979'
$ws.Range('E36').Value = '  -4.04%  '
$ws.Range('D37').Value = '3.055'
$ws.Range('D38').Value = '2.160'
$ws.Range('E38').Value = '  -3.67%  '
$ws.Range('D39').Value = '0.01630'
$ws.Range('E39').Value = '  -7.28%  '
$ws.Range('D40').Value = '0.4421'
$ws.Range('E40').Value = '  -7.62%  '
$ws.Range('D41').Value = '0.8554'
$ws.Range('E41').Value = '  -4.03%  '
$ws.Range('D42').Value = '104.77'
$ws.Range('E42').Value = '  -4.54%  '
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('D44').Value = '5.714'
$ws.Range('E44').Value = '  -3.53%  '
$ws.Range('D45').Value = '6.931'
$ws.Range('E45').Value = '  -9.71%  '
$ws.Range('D46').Value = '8.623'
$ws.Range('E46').Value = '  -3.60%  '
$ws.Range('D47').Value = '0.05557'
$ws.Range('E47').Value = '  -5.00%  '
$ws.Range('D48').Value = '58.36'
$ws.Range('E48').Value = '  -2.07%  '
$ws.Range('D49').Value = '33.23'
$ws.Range('E49').Value = '  -4.58%  '
$ws.Range('D50').Value = '0.1170'
$ws.Range('E50').Value = '  -5.13%  '
$ws.Range('D51').Value = '0.8560'
$ws.Range('E51').Value = '  -4.34%  '
